# ajuste migracion ejecutivos vacios
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New NIT values for the previously-empty "ejecutivos" rows (12-21),
# formatted the same way the original A2:A11 block was before this edit
# (font/border/right-aligned+wrap, no fill).
$newNits = @(900506434,900654612,830032420,900715424,830006055,900947384,830125131,830121765,800142993,830074497)

$startRow = 12
for ($i = 0; $i -lt $newNits.Length; $i++) {
    $row = $startRow + $i
    $src = $ws.Range("A11")
    $src.Copy()
    $dst = $ws.Cells.Item($row, 1)
    $dst.PasteSpecial(-4122) # xlPasteFormats
    $dst.Value = $newNits[$i]
    $ws.Rows.Item($row).RowHeight = 16
}
$excel.CutCopyMode = $false

# Highlight the original (already migrated) rows A2:A11 with a yellow fill
# so they stand out from the newly appended rows.
$ws.Range("A2:A11").Interior.Color = 65535  # RGB(255,255,0)

# Update the visible range / selection to cover the full A2:A21 data block
# and scroll so row 7 is at the top, matching the new working view.
$lastRow = $startRow + $newNits.Length - 1
$ws.Range("A2:A" + $lastRow).Select()
$excel.ActiveWindow.ScrollRow = 7
